$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after the header (rows 2-3), shifting existing data down.
$ws.Rows("2:3").Insert()
$ws.Range("A2:C3").ClearFormats()

$ws.Cells.Item(2, 1).Value = -7.487135410308838
$ws.Cells.Item(2, 2).Value = 2.468842267990112
$ws.Cells.Item(2, 3).Value = -10.7623405456543

$ws.Cells.Item(3, 1).Value = 5.30696439743042
$ws.Cells.Item(3, 2).Value = -7.900328636169434
$ws.Cells.Item(3, 3).Value = -8.124782562255859

# Append eight new rows of data at the bottom (rows 24-31).
$newRows = @(
    @(-7.852145671844482, 6.264562606811523, -14.15214920043945),
    @(-4.107211589813232, -2.712513208389282, -20.56048202514648),
    @(-8.579601287841797, -15.95286655426025, -10.86754608154297),
    @(0.7219026684761047, 2.51579213142395, 17.58181953430176),
    @(-6.065989017486572, 16.71288681030273, -1.314104557037354),
    @(-11.12415027618408, -78.97219848632812, 36.16990280151367),
    @(-0.5275765657424927, 10.35773277282715, -23.80691909790039),
    @(-18.0826530456543, -1.298346519470215, -14.97142791748047)
)

$startRow = 24
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
}
